$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.981.12'
$ws.Range("E2").Value = '  +0.55%  '

$ws.Range("D3").Value = '1.639.62'
$ws.Range("E3").Value = '  -0.04%  '

$ws.Range("E4").Value = '  -0.54%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '217.97'
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = '  +0.14%  '

$ws.Range("E6").Value = '  +1.86%  '

$ws.Range("E7").Value = '  -0.55%  '

$ws.Range("E8").Value = '  +1.55%  '

$ws.Range("E9").Value = '  +0.35%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.98'
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = '  +3.93%  '

$ws.Range("E11").Value = '  -0.01%  '

$ws.Range("D12").Value = '1.869.72'
$ws.Range("E12").Value = '  +0.01%  '

$ws.Range("D13").Value = '1.633.61'
$ws.Range("E13").Value = '  -0.46%  '

$ws.Range("E14").Value = '  -0.81%  '

$ws.Range("E15").Value = '  +1.34%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '67.14'
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = '  +3.03%  '

$ws.Range("D17").Value = '26.974.53'
$ws.Range("E17").Value = '  +0.51%  '

$ws.Range("E18").Value = '  +0.36%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '219.59'
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = '  +2.13%  '

$ws.Range("E21").Value = '  +2.75%  '

$ws.Range("E22").Value = '  +1.36%  '

$ws.Range("E23").Value = '  +1.84%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.20'
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = '  +0.21%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '147.30'
$ws.Range("D25").NumberFormat = "General"

$ws.Range("E26").Value = '  -0.59%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.33'
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = '  +2.00%  '

$ws.Range("E28").Value = '  +0.97%  '

$ws.Range("E29").Value = '  +0.12%  '

$ws.Range("E30").Value = '  -0.72%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.19'
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = '  -0.33%  '

$ws.Range("E32").Value = '  -0.60%  '

$ws.Range("E33").Value = '  +0.58%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.57'
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = '  +1.19%  '

$ws.Range("D35").Value = '1.269.87'
$ws.Range("E35").Value = '  -0.22%  '

$ws.Range("E36").Value = '  -0.13%  '

$ws.Range("E37").Value = '  +2.85%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.540'
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = '  +2.09%  '

$ws.Range("E39").Value = '  +2.66%  '

$ws.Range("E40").Value = '  -0.51%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.808'
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = '  +0.59%  '

$ws.Range("E42").Value = '  +0.73%  '

$ws.Range("D43").Value = '1.780.42'
$ws.Range("E43").Value = '  +0.01%  '

$ws.Range("E44").Value = '  +2.64%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '62.23'
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = '  +2.19%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '92.34'
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = '  -0.08%  '

$ws.Range("E47").Value = '  +1.81%  '

$ws.Range("D48").Value = '0.0₆0103'
$ws.Range("E48").Value = '  -2.26%  '

$ws.Range("E49").Value = '  -0.63%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.65'
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = '  +1.24%  '

$ws.Range("E51").Value = '  -0.37%  '
